# Reorder the weekly price rows (rows 2-9) by ascending date.
# Only columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), Q (Unidad de comercializacion) and
# S (Precio $/Kg) differ between rows; the remaining columns are identical
# across rows and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ D = 44176; M = 250; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 },
    @{ D = 44208; M = 210; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 },
    @{ D = 44351; M = 300; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos empedrada"; S = 714 },
    @{ D = 44309; M = 300; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 },
    @{ D = 44397; M = 60;  N = 11000; O = 11000; P = 11000; Q = "`$/caja 14 kilos";           S = 786 },
    @{ D = 44400; M = 100; N = 10000; O = 10000; P = 10000; Q = "`$/caja 14 kilos";           S = 714 },
    @{ D = 44491; M = 180; N = 9000;  O = 9000;  P = 9000;  Q = "`$/caja 14 kilos empedrada"; S = 643 },
    @{ D = 44162; M = 120; N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 14 kilos empedrada"; S = 500 }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $row = $i + 2
    $data = $rows[$i]

    $ws.Range("D$row").Value = $data.D
    $ws.Range("M$row").Value = $data.M
    $ws.Range("N$row").Value = $data.N
    $ws.Range("O$row").Value = $data.O
    $ws.Range("P$row").Value = $data.P
    $ws.Range("Q$row").Value = $data.Q
    $ws.Range("S$row").Value = $data.S
}
